# Fruta / hortaliza, semanal
# Insert a new weekly record at row 16 (shifting the existing rows 16-65 down
# to 17-66) and populate it with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 16..65 down to 17..66, inserting a fresh row 16.
$ws.Rows(16).Insert()

# Populate the new row 16 with the new weekly record.
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44648
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = "Frutos de pepita"
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = "Membrillo"
$ws.Range("K16").Value = "Champion"
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = 260000
$ws.Range("O16").Value = 280000
$ws.Range("P16").Value = 272000
$ws.Range("Q16").Value = '$/bins (450 kilos)'
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 604
$ws.Range("T16").Value = 450
